# Cenarios de busca por lupa sucesso e falha e busca por home sucesso criados.
#
# This script updates the MassaDeDados.xlsx test-data workbook:
#  - "Cadastro" sheet: refresh the throwaway username used for the cadastro scenario
#  - "PesquisaBarra" sheet: change the searched product / expected-message pair used by
#    the "busca por lupa" (magnifying-glass search) success scenario
#  - "PesquisaCat" sheet: change the expected assert messages used by the "busca por home"
#    success scenario, and make it the active/selected sheet (last worked-on scenario)

$wb = $excel.ActiveWorkbook

$wsCadastro      = $wb.Worksheets.Item("Cadastro")
$wsPesquisaCat   = $wb.Worksheets.Item("PesquisaCat")
$wsPesquisaBarra = $wb.Worksheets.Item("PesquisaBarra")

# --- PesquisaBarra: novo produto buscado / mensagem esperada ----------------
# B2 keeps its plain style; just overwrite the value.
$wsPesquisaBarra.Activate()
$wsPesquisaBarra.Range("B2").Value = "HP H2310 IN-EAR HEADSET"
# E2 uses a quote-prefixed ("text") style in the original file; a leading
# apostrophe preserves that quote-prefix formatting instead of resetting it.
$wsPesquisaBarra.Range("E2").Value = "'HP H2310 IN-EAR HEADSET"
$wsPesquisaBarra.Range("B2").Select()

# --- Cadastro: novo usuario de teste ---------------------------------------
$wsCadastro.Activate()
$wsCadastro.Range("B2").Value = "lucascarvalh23"
$wsCadastro.Range("B2").Select()

# --- PesquisaCat: mensagens de assert da busca por home ---------------------
$wsPesquisaCat.Activate()
$wsPesquisaCat.Range("E2").Value = "Pesquisa efetuada com sucesso!"
$wsPesquisaCat.Range("E3").Value = "HP ROAR MINI WIRELESS SPEAKER"
$wsPesquisaCat.Range("E3").Select()

# Leave PesquisaCat as the active/selected tab, matching the last scenario worked on.
$wsPesquisaCat.Activate()
